$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new "Concepto" column at D, shifting De/Importe/USD Prom/Dolares
#    one column to the right (D->E, E->F, F->G, G->H). Formulas auto-adjust.
$ws.Columns("D:D").Insert()

$ws.Range("D2").Value = "Concepto"
$ws.Range("D3").Value = "Abono Anual 2024 "
$ws.Range("D4").Value = "Abono Anual 2024 "
$ws.Range("D5").Value = "Abono Anual 2024 "

# Column D picked up column C's per-row formatting on insert; row 5 needs to
# match the "date family" style (like B3) instead of C5's plain style.
$ws.Range("B3").Copy()
$ws.Range("D5").PasteSpecial(-4122)

# 2. New row 6: "Silvia Barros Reyes" / "Extra x Dif pago Tarjeta" entry.
$ws.Range("B6").Value = 45457
$ws.Range("C6").Value = "Transferencia Mariano"
$ws.Range("D6").Value = "Abono Anual 2024 "
$ws.Range("E6").Value = "Silvia Barros Reyes"
$ws.Range("F6").Value = 33000
$ws.Range("G6").Value = 2010.44
$ws.Range("H6").Formula = "=+F6/G6"
$ws.Range("I6").Value = "Extra x Dif pago Tarjeta"

# Copy row 5's per-column formatting down into row 6 (date/text/currency/usd
# families), then the highlight fill below is layered on top of it.
$ws.Range("B5:H5").Copy()
$ws.Range("B6:H6").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# 3. Highlight the whole data block (new "Abono Anual 2024" rows) in yellow.
$ws.Range("B3:H6").Interior.Color = 65535

# 4. Threaded comment flagging the real USD value for the new row.
$ws.Range("G6").AddCommentThreaded("Valor real USD 1280") | Out-Null

# 5. Column widths to fit the new layout.
$ws.Columns("B").ColumnWidth = 10.42578125
$ws.Columns("D").ColumnWidth = 21.28515625
$ws.Columns("E").ColumnWidth = 23.140625
$ws.Columns("I").ColumnWidth = 21
$ws.Columns("J").ColumnWidth = 13.5703125

$ws.Range("C15").Select()
